# Update "想去人数" (want-to-go count) figures across the four sheets, and
# add the new "北京·Paradox Live·[灯光渐强]" row to the 本地生活 (Local
# life) sheet, matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 440
$ws.Range("F4").Value = 124
$ws.Range("F6").Value = 507
$ws.Range("F7").Value = 40
$ws.Range("F9").Value = 255
$ws.Range("F10").Value = 369
$ws.Range("F12").Value = 551
$ws.Range("F13").Value = 719
$ws.Range("F14").Value = 1475
$ws.Range("F15").Value = 1475
$ws.Range("F16").Value = 875
$ws.Range("F18").Value = 1333
$ws.Range("F20").Value = 238
$ws.Range("F23").Value = 89
$ws.Range("F24").Value = 6374
$ws.Range("F25").Value = 4699
$ws.Range("F26").Value = 121
$ws.Range("F27").Value = 485
$ws.Range("F28").Value = 145
$ws.Range("F29").Value = 61
$ws.Range("F32").Value = 1228
$ws.Range("F33").Value = 177
$ws.Range("F34").Value = 27
$ws.Range("F35").Value = 576
$ws.Range("F38").Value = 209
$ws.Range("F40").Value = 133
$ws.Range("F42").Value = 82
$ws.Range("F44").Value = 7

# ---------------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 1

# ---------------------------------------------------------------------
# Sheet: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 167

# New row 5: copy row 4's column-A formatting (bold / bordered / centred
# index style) onto the new index cell, then overwrite with the new data.
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 4

# Force text entry on the date cell so the "YYYY-MM-DD" string isn't
# auto-converted into a date serial number (matches the existing B column
# cells, which are stored as plain text); ClearFormats afterwards so the
# cell keeps the plain/default style like its siblings (no explicit "s"
# attribute), same as the rest of column B.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2024-08-28"
$ws.Range("B5").ClearFormats()

$ws.Range("C5").Value = "北京·Paradox Live·[灯光渐强] "
$ws.Range("D5").Value = "王府井地铁站F1东口步行120米 北京王府井喜悦购物中心"
$ws.Range("E5").Value = "2024.08.28 00:00-10.11 23:59"
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=91230"
$ws.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202408/WV8PxG321724639038452.jpeg"

# ---------------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 440
$ws.Range("F7").Value = 167
$ws.Range("F9").Value = 507
$ws.Range("F10").Value = 40
$ws.Range("F12").Value = 255
$ws.Range("F14").Value = 369
$ws.Range("F16").Value = 551
$ws.Range("F17").Value = 719
$ws.Range("F18").Value = 1475
$ws.Range("F19").Value = 1475
$ws.Range("F20").Value = 875
$ws.Range("F22").Value = 1333
$ws.Range("F24").Value = 238
$ws.Range("F27").Value = 89
$ws.Range("F30").Value = 6374
$ws.Range("F31").Value = 4699
$ws.Range("F32").Value = 121
$ws.Range("F34").Value = 1228
$ws.Range("F35").Value = 177
$ws.Range("F36").Value = 27
$ws.Range("F38").Value = 576
$ws.Range("F44").Value = 210
$ws.Range("F45").Value = 133
$ws.Range("F47").Value = 82
$ws.Range("F50").Value = 7
